$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Keep gridlines visible on the original sheet (engine defaults to hidden
# once the sheetView is rewritten, so force it back on)
$ws1.Application.ActiveWindow.DisplayGridlines = $true

# Copy the finalized daily-data table (header row + 30 days) into a new
# sheet, preserving formatting
$ws1.Range("A9:K39").Copy() | Out-Null

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet1"
$ws2.Range("A1").PasteSpecial() | Out-Null
$ws2.Range("A1:K31").Select() | Out-Null

# Scroll/select the source sheet back to the table region
$ws1.Range("A9:K39").Select() | Out-Null

# The new recap sheet becomes the active tab
$ws2.Activate()
